$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.225.35"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "2.306.62"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.510"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.500"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.96%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.73%  "
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.64%  "
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "2.656.38"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").Value = "2.312.17"
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.813"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").Value = "43.002.40"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.26%  "
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "237.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -3.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  +3.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.28%  "
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("E38").Value = "  -2.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  -3.42%  "
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("E43").Value = "  -4.28%  "
$ws.Range("D44").Value = "1.989.06"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0282"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.08%  "
$ws.Range("E48").Value = "  -2.59%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.527.66"
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.18%  "
